# Add a new "QA Engineer" job row to the jobs sheet and fix punctuation on
# the "Admin" row's description (trailing period added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: QA Engineer in Chennai.
# (Operations ordered so new shared strings are appended in the same order
# seen in the target workbook: QA Engineer, Chennai, Test and deliver...)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "QA Engineer"
$ws.Range("C6").Value = "Chennai"
$ws.Range("D6").Value = "Test and deliver quality products."

# Existing Admin row (row 5) description gets a trailing period added.
$ws.Range("D5").Value = "Manage daily tasks and office tasks."

# Match the author's final selection/view state.
$ws.Range("D6").Select()
